$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.558.54'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.479.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9766'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '279.37'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3663'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.75%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.97'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.062'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06657'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.507'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.61%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9774'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001028'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.481.24'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05929'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.50'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.475'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.28%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.248'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.604.37'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.41'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.148'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -7.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.27'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.637.12'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.75'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.993'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8138'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08038'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.548'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.232'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +8.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05832'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.725'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.66%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.751'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.73%  '
$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9767'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02048'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.46'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1888'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5296'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.13'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.82'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5201'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.799'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06464'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9934'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.82%  '
